$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete row 8 entirely (q1411379 moves up into what becomes row 7's data set,
# other rows shift as the IDs in rows 5-7 change to reflect the row that was removed)
$ws.Rows.Item(8).Delete()

# New values for row 2 (c1206235)
$ws.Range("B2").Value = 0.84143
$ws.Range("C2").Value = 0.06779
$ws.Range("D2").Value = 0.59297
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 5
$ws.Range("G2").Value = 7

# New values for row 3 (c1242115)
$ws.Range("B3").Value = 1.6699
$ws.Range("C3").Value = 0.62208
$ws.Range("D3").Value = 0.8920400000000001
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 5
$ws.Range("G3").Value = 9

# New values for row 4 (c1243957)
$ws.Range("B4").Value = 1.0695
$ws.Range("C4").Value = 0.67005
$ws.Range("D4").Value = 0.82267
$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 8

# Row 5: ID changes from q0328135_previewuser to q0762379, plus new values
$ws.Range("A5").Value = "q0762379"
$ws.Range("B5").Value = 0.9403899999999999
$ws.Range("C5").Value = 0.71472
$ws.Range("D5").Value = 0.3838
$ws.Range("E5").Value = 1

# Row 6: ID changes from q0762379 to q1371623, plus new values
$ws.Range("A6").Value = "q1371623"
$ws.Range("B6").Value = 2.26283
$ws.Range("C6").Value = 0.5778
$ws.Range("D6").Value = 0.6782899999999999
$ws.Range("F6").Value = 5

# Row 7: ID changes from q1371623 to q1411379, plus new values
$ws.Range("A7").Value = "q1411379"
$ws.Range("B7").Value = 1.24665
$ws.Range("C7").Value = 0.98333
$ws.Range("D7").Value = 0.63576
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 7
